# Video09b / 5 case studies title slide update:
# Prepend "Video09b - " (as separate word/space/dash runs, matching the
# existing per-word run layout already used for the rest of the title)
# to the title placeholder on the first slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Insert in reverse order so each new word lands before the previous one,
# ending up as four distinct runs: "Video09b", " ", "-", " ".
# ($null = ...) swallows the COM return value so it doesn't leak into the
# script's output stream.
$null = $tr.InsertBefore(" ")
$null = $tr.InsertBefore("-")
$null = $tr.InsertBefore(" ")
$null = $tr.InsertBefore("Video09b")
